$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 122, shifting existing rows 122-125 down to 123-126.
$ws.Rows("122:122").Insert()

# The date column (D) uses a custom date style; copy that single cell's
# formatting from the row below (which used to be row 122) onto the new cell.
$ws.Range("D123").Copy()
$ws.Range("D122").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A122").Value = 4
$ws.Range("B122").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C122").Value = "Los Lagos"
$ws.Range("D122").Value = 44448
$ws.Range("E122").Value = 10
$ws.Range("F122").Value = 100112017
$ws.Range("G122").Value = "Apio"
$ws.Range("H122").Value = "Americana (o)"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 25
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = 12000
$ws.Range("N122").Value = "`$/docena de matas"
$ws.Range("O122").Value = "Región de Coquimbo"
$ws.Range("P122").Value = 2000
$ws.Range("Q122").Value = 6
$ws.Range("R122").Value = "Hortaliza"
